$wb = $excel.ActiveWorkbook

# ===== Sheet 1 (sheet1) =====
$ws = $wb.Worksheets.Item(1)

# Update cell values
$ws.Range("A1").Value = "File Name"
$ws.Range("B1").Value = "zh-cn"
$ws.Range("C1").Value = "de-de"
$ws.Range("D1").Value = "Latest Handoff Date"
$ws.Range("A2").Value = "ac08f94f-0e8f-46ad-8098-6c5698180711.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "2016-20-19 02:20:07"
$ws.Range("A3").Value = "be2289b2-8feb-4877-b051-3dfe47b31b49.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "2016-17-19 02:17:15"
$ws.Range("A4").Value = "dc985425-1c7c-4fd9-89fd-e6724d4c27d2.md"
$ws.Range("B4").Value = "Handed back: in sync with en-US"
$ws.Range("C4").Value = "Handed back: in sync with en-US"
$ws.Range("D4").Value = "2016-20-19 02:20:07"
$ws.Range("A5").Value = "e988cb26-8b22-4b05-b024-22b86466f393.md"
$ws.Range("B5").Value = "Handed back: in sync with en-US"
$ws.Range("C5").Value = "Handed back: in sync with en-US"
$ws.Range("D5").Value = "2016-17-19 02:17:15"
$ws.Range("A6").Value = "a8fe6dd2-032f-4017-b24f-9e9414829f06.md"
$ws.Range("B6").Value = "Handback transform failed"
$ws.Range("C6").Value = "Handback transform failed"
$ws.Range("D6").Value = "2016-20-19 02:20:07"
$ws.Range("A7").Value = "05bdf957-5021-4c9f-bbc8-fe4e026e9c96.md"
$ws.Range("B7").Value = "In Translation"
$ws.Range("C7").Value = "In Translation"
$ws.Range("D7").Value = "2016-16-19 02:16:16"
$ws.Range("A8").Value = "73564754-f334-4035-a774-30b6bcac73ee.md"
$ws.Range("B8").Value = "Ready for handoff"
$ws.Range("C8").Value = "Ready for handoff"
$ws.Range("D8").Value = "2016-20-19 02:20:07"

# Rebuild hyperlinks to match final row order/content
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7987b92ebd0d32149992965e10f2143cdfc116df/e2e/ac08f94f-0e8f-46ad-8098-6c5698180711.md", "", "", "ac08f94f-0e8f-46ad-8098-6c5698180711.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d03787431f5a3d598fef0f50f326d128bb146a3d/e2e/be2289b2-8feb-4877-b051-3dfe47b31b49.md", "", "", "be2289b2-8feb-4877-b051-3dfe47b31b49.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/994cba010e790f58810ae1a218b0cc121843b6bc/e2e/dc985425-1c7c-4fd9-89fd-e6724d4c27d2.md", "", "", "dc985425-1c7c-4fd9-89fd-e6724d4c27d2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/d03787431f5a3d598fef0f50f326d128bb146a3d/e2e/e988cb26-8b22-4b05-b024-22b86466f393.md", "", "", "e988cb26-8b22-4b05-b024-22b86466f393.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/3262cec2e5534eb1fb51addfd3cfed290b4fcb0e/e2e/05bdf957-5021-4c9f-bbc8-fe4e026e9c96.md", "", "", "a8fe6dd2-032f-4017-b24f-9e9414829f06.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/5ec76fc78b7923d3c111025329fa9d307811021d/e2e/73564754-f334-4035-a774-30b6bcac73ee.md", "", "", "05bdf957-5021-4c9f-bbc8-fe4e026e9c96.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/9b6d09d209c958f904d518b6d5614b54e38cf1fd/e2e/a8fe6dd2-032f-4017-b24f-9e9414829f06.md", "", "", "73564754-f334-4035-a774-30b6bcac73ee.md") | Out-Null

# ===== Sheet 2 (sheet2) =====
$ws = $wb.Worksheets.Item(2)

# Update cell values
$ws.Range("A1").Value = "Source File Name"
$ws.Range("B1").Value = "File Extension"
$ws.Range("C1").Value = "Status"
$ws.Range("D1").Value = "Latest Handoff File"
$ws.Range("E1").Value = "Latest Handoff Datetime"
$ws.Range("F1").Value = "Latest Target File"
$ws.Range("G1").Value = "Latest Handback File"
$ws.Range("H1").Value = "Latest Handback DateTime"
$ws.Range("I1").Value = "Handoff Reason"
$ws.Range("J1").Value = "Dependency From"
$ws.Range("K1").Value = "Error Detail"
$ws.Range("A2").Value = "ac08f94f-0e8f-46ad-8098-6c5698180711.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "ac08f94f-0e8f-46ad-8098-6c5698180711.235d34e039b5bb95c22b481201035261684078a4.zh-cn.xlf"
$ws.Range("E2").Value = "2016-03-19 02:20:02"
$ws.Range("F2").Value = "ac08f94f-0e8f-46ad-8098-6c5698180711.md"
$ws.Range("G2").Value = "ac08f94f-0e8f-46ad-8098-6c5698180711.235d34e039b5bb95c22b481201035261684078a4.zh-cn.xlf"
$ws.Range("H2").Value = "2016-03-19 02:19:29"
$ws.Range("I2").Value = "Include"
$ws.Range("A3").Value = "be2289b2-8feb-4877-b051-3dfe47b31b49.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "be2289b2-8feb-4877-b051-3dfe47b31b49.4a23415bd45c3cd95c9d4a0e51ba1905a5991e31.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-19 02:17:12"
$ws.Range("F3").Value = "be2289b2-8feb-4877-b051-3dfe47b31b49.md"
$ws.Range("G3").Value = "be2289b2-8feb-4877-b051-3dfe47b31b49.4a23415bd45c3cd95c9d4a0e51ba1905a5991e31.zh-cn.xlf"
$ws.Range("H3").Value = "2016-03-19 02:17:30"
$ws.Range("I3").Value = "Include"
$ws.Range("A4").Value = "dc985425-1c7c-4fd9-89fd-e6724d4c27d2.md"
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "Handed back: in sync with en-US"
$ws.Range("D4").Value = "dc985425-1c7c-4fd9-89fd-e6724d4c27d2.562879b4d3044da5adb5f238ed954a3eb20d8851.zh-cn.xlf"
$ws.Range("E4").Value = "2016-03-19 02:20:02"
$ws.Range("F4").Value = "dc985425-1c7c-4fd9-89fd-e6724d4c27d2.md"
$ws.Range("G4").Value = "dc985425-1c7c-4fd9-89fd-e6724d4c27d2.562879b4d3044da5adb5f238ed954a3eb20d8851.zh-cn.xlf"
$ws.Range("H4").Value = "2016-03-19 02:19:29"
$ws.Range("I4").Value = "Include"
$ws.Range("A5").Value = "e988cb26-8b22-4b05-b024-22b86466f393.md"
$ws.Range("B5").Value = ".md"
$ws.Range("C5").Value = "Handed back: in sync with en-US"
$ws.Range("D5").Value = "e988cb26-8b22-4b05-b024-22b86466f393.217cbbfe5e36df50d96a70fd14c971ac31104000.zh-cn.xlf"
$ws.Range("E5").Value = "2016-03-19 02:17:12"
$ws.Range("F5").Value = "e988cb26-8b22-4b05-b024-22b86466f393.md"
$ws.Range("G5").Value = "e988cb26-8b22-4b05-b024-22b86466f393.217cbbfe5e36df50d96a70fd14c971ac31104000.zh-cn.xlf"
$ws.Range("H5").Value = "2016-03-19 02:17:30"
$ws.Range("I5").Value = "Include"
$ws.Range("A6").Value = "a8fe6dd2-032f-4017-b24f-9e9414829f06.md"
$ws.Range("B6").Value = ".md"
$ws.Range("C6").Value = "Handback transform failed"
$ws.Range("D6").Value = "a8fe6dd2-032f-4017-b24f-9e9414829f06.4b0e106ca6a6d67c239a240617393aacd849fcaa.zh-cn.xlf"
$ws.Range("E6").Value = "2016-03-19 02:20:02"
$ws.Range("H6").Value = "0001-01-01 00:00:00"
$ws.Range("I6").Value = "Include"
$ws.Range("K6").Value = "The handback type mt is not match with handoff type ht."
$ws.Range("A7").Value = "05bdf957-5021-4c9f-bbc8-fe4e026e9c96.md"
$ws.Range("B7").Value = ".md"
$ws.Range("C7").Value = "In Translation"
$ws.Range("D7").Value = "05bdf957-5021-4c9f-bbc8-fe4e026e9c96.0d47ba9f2dce54aba0760192d3c7c3ca0bb473a8.zh-cn.xlf"
$ws.Range("E7").Value = "2016-03-19 02:16:13"
$ws.Range("H7").Value = "0001-01-01 00:00:00"
$ws.Range("I7").Value = "Include"
$ws.Range("A8").Value = "73564754-f334-4035-a774-30b6bcac73ee.md"
$ws.Range("B8").Value = ".md"
$ws.Range("C8").Value = "Ready for handoff"
$ws.Range("D8").Value = "73564754-f334-4035-a774-30b6bcac73ee.a042b5a7eea815cab0eb25d01f00914718577d55.zh-cn.xlf"
$ws.Range("E8").Value = "2016-03-19 02:20:02"
$ws.Range("H8").Value = "0001-01-01 00:00:00"
$ws.Range("I8").Value = "Include"

# Rebuild hyperlinks to match final row order/content
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7987b92ebd0d32149992965e10f2143cdfc116df/e2e/ac08f94f-0e8f-46ad-8098-6c5698180711.md", "", "", "ac08f94f-0e8f-46ad-8098-6c5698180711.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/7987b92ebd0d32149992965e10f2143cdfc116df/e2e/ac08f94f-0e8f-46ad-8098-6c5698180711.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4738c07cbedf7b7a0c0c478a90de59463c80927b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ac08f94f-0e8f-46ad-8098-6c5698180711.235d34e039b5bb95c22b481201035261684078a4.zh-cn.xlf", "", "", "ac08f94f-0e8f-46ad-8098-6c5698180711.235d34e039b5bb95c22b481201035261684078a4.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/48ab36ef6d69f76e8d1de60f981bbac7fb820ff2/e2e/ac08f94f-0e8f-46ad-8098-6c5698180711.md", "", "", "ac08f94f-0e8f-46ad-8098-6c5698180711.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/073d4fb53bf6803d062ecccba3289ef0535f7008/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/mt/ac08f94f-0e8f-46ad-8098-6c5698180711.235d34e039b5bb95c22b481201035261684078a4.zh-cn.xlf", "", "", "ac08f94f-0e8f-46ad-8098-6c5698180711.235d34e039b5bb95c22b481201035261684078a4.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d03787431f5a3d598fef0f50f326d128bb146a3d/e2e/be2289b2-8feb-4877-b051-3dfe47b31b49.md", "", "", "be2289b2-8feb-4877-b051-3dfe47b31b49.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/d03787431f5a3d598fef0f50f326d128bb146a3d/e2e/be2289b2-8feb-4877-b051-3dfe47b31b49.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/97b7ad6ce55f104cdc84d036f4bddf08477f0159/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/be2289b2-8feb-4877-b051-3dfe47b31b49.4a23415bd45c3cd95c9d4a0e51ba1905a5991e31.zh-cn.xlf", "", "", "be2289b2-8feb-4877-b051-3dfe47b31b49.4a23415bd45c3cd95c9d4a0e51ba1905a5991e31.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/fae95c2285788959fc0222af960f618587ca0eac/e2e/be2289b2-8feb-4877-b051-3dfe47b31b49.md", "", "", "be2289b2-8feb-4877-b051-3dfe47b31b49.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b559cba1dc73597d739b343c2ba99e3ecc0144db/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/be2289b2-8feb-4877-b051-3dfe47b31b49.4a23415bd45c3cd95c9d4a0e51ba1905a5991e31.zh-cn.xlf", "", "", "be2289b2-8feb-4877-b051-3dfe47b31b49.4a23415bd45c3cd95c9d4a0e51ba1905a5991e31.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/994cba010e790f58810ae1a218b0cc121843b6bc/e2e/dc985425-1c7c-4fd9-89fd-e6724d4c27d2.md", "", "", "dc985425-1c7c-4fd9-89fd-e6724d4c27d2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/994cba010e790f58810ae1a218b0cc121843b6bc/e2e/dc985425-1c7c-4fd9-89fd-e6724d4c27d2.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4738c07cbedf7b7a0c0c478a90de59463c80927b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/dc985425-1c7c-4fd9-89fd-e6724d4c27d2.562879b4d3044da5adb5f238ed954a3eb20d8851.zh-cn.xlf", "", "", "dc985425-1c7c-4fd9-89fd-e6724d4c27d2.562879b4d3044da5adb5f238ed954a3eb20d8851.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/48ab36ef6d69f76e8d1de60f981bbac7fb820ff2/e2e/dc985425-1c7c-4fd9-89fd-e6724d4c27d2.md", "", "", "dc985425-1c7c-4fd9-89fd-e6724d4c27d2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/073d4fb53bf6803d062ecccba3289ef0535f7008/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/mt/dc985425-1c7c-4fd9-89fd-e6724d4c27d2.562879b4d3044da5adb5f238ed954a3eb20d8851.zh-cn.xlf", "", "", "dc985425-1c7c-4fd9-89fd-e6724d4c27d2.562879b4d3044da5adb5f238ed954a3eb20d8851.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/d03787431f5a3d598fef0f50f326d128bb146a3d/e2e/e988cb26-8b22-4b05-b024-22b86466f393.md", "", "", "e988cb26-8b22-4b05-b024-22b86466f393.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/d03787431f5a3d598fef0f50f326d128bb146a3d/e2e/e988cb26-8b22-4b05-b024-22b86466f393.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/97b7ad6ce55f104cdc84d036f4bddf08477f0159/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/e988cb26-8b22-4b05-b024-22b86466f393.217cbbfe5e36df50d96a70fd14c971ac31104000.zh-cn.xlf", "", "", "e988cb26-8b22-4b05-b024-22b86466f393.217cbbfe5e36df50d96a70fd14c971ac31104000.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/fae95c2285788959fc0222af960f618587ca0eac/e2e/e988cb26-8b22-4b05-b024-22b86466f393.md", "", "", "e988cb26-8b22-4b05-b024-22b86466f393.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b559cba1dc73597d739b343c2ba99e3ecc0144db/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e988cb26-8b22-4b05-b024-22b86466f393.217cbbfe5e36df50d96a70fd14c971ac31104000.zh-cn.xlf", "", "", "e988cb26-8b22-4b05-b024-22b86466f393.217cbbfe5e36df50d96a70fd14c971ac31104000.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/3262cec2e5534eb1fb51addfd3cfed290b4fcb0e/e2e/05bdf957-5021-4c9f-bbc8-fe4e026e9c96.md", "", "", "a8fe6dd2-032f-4017-b24f-9e9414829f06.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B6"), "https://github.com/OpenLocalizationTest/oltest/blob/3262cec2e5534eb1fb51addfd3cfed290b4fcb0e/e2e/05bdf957-5021-4c9f-bbc8-fe4e026e9c96.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/15f7fa9289cee50a4f37d3f67c10feaadbbe2b1d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/05bdf957-5021-4c9f-bbc8-fe4e026e9c96.0d47ba9f2dce54aba0760192d3c7c3ca0bb473a8.zh-cn.xlf", "", "", "a8fe6dd2-032f-4017-b24f-9e9414829f06.4b0e106ca6a6d67c239a240617393aacd849fcaa.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/5ec76fc78b7923d3c111025329fa9d307811021d/e2e/73564754-f334-4035-a774-30b6bcac73ee.md", "", "", "05bdf957-5021-4c9f-bbc8-fe4e026e9c96.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B7"), "https://github.com/OpenLocalizationTest/oltest/blob/5ec76fc78b7923d3c111025329fa9d307811021d/e2e/73564754-f334-4035-a774-30b6bcac73ee.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4738c07cbedf7b7a0c0c478a90de59463c80927b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/73564754-f334-4035-a774-30b6bcac73ee.a042b5a7eea815cab0eb25d01f00914718577d55.zh-cn.xlf", "", "", "05bdf957-5021-4c9f-bbc8-fe4e026e9c96.0d47ba9f2dce54aba0760192d3c7c3ca0bb473a8.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/9b6d09d209c958f904d518b6d5614b54e38cf1fd/e2e/a8fe6dd2-032f-4017-b24f-9e9414829f06.md", "", "", "73564754-f334-4035-a774-30b6bcac73ee.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B8"), "https://github.com/OpenLocalizationTest/oltest/blob/9b6d09d209c958f904d518b6d5614b54e38cf1fd/e2e/a8fe6dd2-032f-4017-b24f-9e9414829f06.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D8"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4738c07cbedf7b7a0c0c478a90de59463c80927b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a8fe6dd2-032f-4017-b24f-9e9414829f06.4b0e106ca6a6d67c239a240617393aacd849fcaa.zh-cn.xlf", "", "", "73564754-f334-4035-a774-30b6bcac73ee.a042b5a7eea815cab0eb25d01f00914718577d55.zh-cn.xlf") | Out-Null

# ===== Sheet 3 (sheet3) =====
$ws = $wb.Worksheets.Item(3)

# Update cell values
$ws.Range("A1").Value = "Source File Name"
$ws.Range("B1").Value = "File Extension"
$ws.Range("C1").Value = "Status"
$ws.Range("D1").Value = "Latest Handoff File"
$ws.Range("E1").Value = "Latest Handoff Datetime"
$ws.Range("F1").Value = "Latest Target File"
$ws.Range("G1").Value = "Latest Handback File"
$ws.Range("H1").Value = "Latest Handback DateTime"
$ws.Range("I1").Value = "Handoff Reason"
$ws.Range("J1").Value = "Dependency From"
$ws.Range("K1").Value = "Error Detail"
$ws.Range("A2").Value = "ac08f94f-0e8f-46ad-8098-6c5698180711.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "ac08f94f-0e8f-46ad-8098-6c5698180711.235d34e039b5bb95c22b481201035261684078a4.de-de.xlf"
$ws.Range("E2").Value = "2016-03-19 02:20:07"
$ws.Range("F2").Value = "ac08f94f-0e8f-46ad-8098-6c5698180711.md"
$ws.Range("G2").Value = "ac08f94f-0e8f-46ad-8098-6c5698180711.235d34e039b5bb95c22b481201035261684078a4.de-de.xlf"
$ws.Range("H2").Value = "2016-03-19 02:19:34"
$ws.Range("I2").Value = "Include"
$ws.Range("A3").Value = "be2289b2-8feb-4877-b051-3dfe47b31b49.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "be2289b2-8feb-4877-b051-3dfe47b31b49.4a23415bd45c3cd95c9d4a0e51ba1905a5991e31.de-de.xlf"
$ws.Range("E3").Value = "2016-03-19 02:17:15"
$ws.Range("F3").Value = "be2289b2-8feb-4877-b051-3dfe47b31b49.md"
$ws.Range("G3").Value = "be2289b2-8feb-4877-b051-3dfe47b31b49.4a23415bd45c3cd95c9d4a0e51ba1905a5991e31.de-de.xlf"
$ws.Range("H3").Value = "2016-03-19 02:17:35"
$ws.Range("I3").Value = "Include"
$ws.Range("A4").Value = "dc985425-1c7c-4fd9-89fd-e6724d4c27d2.md"
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "Handed back: in sync with en-US"
$ws.Range("D4").Value = "dc985425-1c7c-4fd9-89fd-e6724d4c27d2.562879b4d3044da5adb5f238ed954a3eb20d8851.de-de.xlf"
$ws.Range("E4").Value = "2016-03-19 02:20:07"
$ws.Range("F4").Value = "dc985425-1c7c-4fd9-89fd-e6724d4c27d2.md"
$ws.Range("G4").Value = "dc985425-1c7c-4fd9-89fd-e6724d4c27d2.562879b4d3044da5adb5f238ed954a3eb20d8851.de-de.xlf"
$ws.Range("H4").Value = "2016-03-19 02:19:34"
$ws.Range("I4").Value = "Include"
$ws.Range("A5").Value = "e988cb26-8b22-4b05-b024-22b86466f393.md"
$ws.Range("B5").Value = ".md"
$ws.Range("C5").Value = "Handed back: in sync with en-US"
$ws.Range("D5").Value = "e988cb26-8b22-4b05-b024-22b86466f393.217cbbfe5e36df50d96a70fd14c971ac31104000.de-de.xlf"
$ws.Range("E5").Value = "2016-03-19 02:17:15"
$ws.Range("F5").Value = "e988cb26-8b22-4b05-b024-22b86466f393.md"
$ws.Range("G5").Value = "e988cb26-8b22-4b05-b024-22b86466f393.217cbbfe5e36df50d96a70fd14c971ac31104000.de-de.xlf"
$ws.Range("H5").Value = "2016-03-19 02:17:35"
$ws.Range("I5").Value = "Include"
$ws.Range("A6").Value = "a8fe6dd2-032f-4017-b24f-9e9414829f06.md"
$ws.Range("B6").Value = ".md"
$ws.Range("C6").Value = "Handback transform failed"
$ws.Range("D6").Value = "a8fe6dd2-032f-4017-b24f-9e9414829f06.4b0e106ca6a6d67c239a240617393aacd849fcaa.de-de.xlf"
$ws.Range("E6").Value = "2016-03-19 02:20:07"
$ws.Range("H6").Value = "0001-01-01 00:00:00"
$ws.Range("I6").Value = "Include"
$ws.Range("K6").Value = "The handback type mt is not match with handoff type ht."
$ws.Range("A7").Value = "05bdf957-5021-4c9f-bbc8-fe4e026e9c96.md"
$ws.Range("B7").Value = ".md"
$ws.Range("C7").Value = "In Translation"
$ws.Range("D7").Value = "05bdf957-5021-4c9f-bbc8-fe4e026e9c96.0d47ba9f2dce54aba0760192d3c7c3ca0bb473a8.de-de.xlf"
$ws.Range("E7").Value = "2016-03-19 02:16:16"
$ws.Range("H7").Value = "0001-01-01 00:00:00"
$ws.Range("I7").Value = "Include"
$ws.Range("A8").Value = "73564754-f334-4035-a774-30b6bcac73ee.md"
$ws.Range("B8").Value = ".md"
$ws.Range("C8").Value = "Ready for handoff"
$ws.Range("D8").Value = "73564754-f334-4035-a774-30b6bcac73ee.a042b5a7eea815cab0eb25d01f00914718577d55.de-de.xlf"
$ws.Range("E8").Value = "2016-03-19 02:20:07"
$ws.Range("H8").Value = "0001-01-01 00:00:00"
$ws.Range("I8").Value = "Include"

# Rebuild hyperlinks to match final row order/content
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7987b92ebd0d32149992965e10f2143cdfc116df/e2e/ac08f94f-0e8f-46ad-8098-6c5698180711.md", "", "", "ac08f94f-0e8f-46ad-8098-6c5698180711.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/7987b92ebd0d32149992965e10f2143cdfc116df/e2e/ac08f94f-0e8f-46ad-8098-6c5698180711.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5c0affbe02d6842e482f38043da5d48ac8a240c0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ac08f94f-0e8f-46ad-8098-6c5698180711.235d34e039b5bb95c22b481201035261684078a4.de-de.xlf", "", "", "ac08f94f-0e8f-46ad-8098-6c5698180711.235d34e039b5bb95c22b481201035261684078a4.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/2e4b41fbbd5b136a25cbe39beedd4977163f7d36/e2e/ac08f94f-0e8f-46ad-8098-6c5698180711.md", "", "", "ac08f94f-0e8f-46ad-8098-6c5698180711.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/018e5cf3694a27eed9e0a56c9b9439f952746b5a/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/mt/ac08f94f-0e8f-46ad-8098-6c5698180711.235d34e039b5bb95c22b481201035261684078a4.de-de.xlf", "", "", "ac08f94f-0e8f-46ad-8098-6c5698180711.235d34e039b5bb95c22b481201035261684078a4.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d03787431f5a3d598fef0f50f326d128bb146a3d/e2e/be2289b2-8feb-4877-b051-3dfe47b31b49.md", "", "", "be2289b2-8feb-4877-b051-3dfe47b31b49.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/d03787431f5a3d598fef0f50f326d128bb146a3d/e2e/be2289b2-8feb-4877-b051-3dfe47b31b49.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/49f891363dce45ef12fb8a60c2bc2d061cbeb1c4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/be2289b2-8feb-4877-b051-3dfe47b31b49.4a23415bd45c3cd95c9d4a0e51ba1905a5991e31.de-de.xlf", "", "", "be2289b2-8feb-4877-b051-3dfe47b31b49.4a23415bd45c3cd95c9d4a0e51ba1905a5991e31.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/243c54956fc46f4693e9a4921faacc08b9574e2d/e2e/be2289b2-8feb-4877-b051-3dfe47b31b49.md", "", "", "be2289b2-8feb-4877-b051-3dfe47b31b49.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2a2c62b0eb4ee1edf095955e40781832d0c9b312/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/be2289b2-8feb-4877-b051-3dfe47b31b49.4a23415bd45c3cd95c9d4a0e51ba1905a5991e31.de-de.xlf", "", "", "be2289b2-8feb-4877-b051-3dfe47b31b49.4a23415bd45c3cd95c9d4a0e51ba1905a5991e31.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/994cba010e790f58810ae1a218b0cc121843b6bc/e2e/dc985425-1c7c-4fd9-89fd-e6724d4c27d2.md", "", "", "dc985425-1c7c-4fd9-89fd-e6724d4c27d2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/994cba010e790f58810ae1a218b0cc121843b6bc/e2e/dc985425-1c7c-4fd9-89fd-e6724d4c27d2.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5c0affbe02d6842e482f38043da5d48ac8a240c0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/dc985425-1c7c-4fd9-89fd-e6724d4c27d2.562879b4d3044da5adb5f238ed954a3eb20d8851.de-de.xlf", "", "", "dc985425-1c7c-4fd9-89fd-e6724d4c27d2.562879b4d3044da5adb5f238ed954a3eb20d8851.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/2e4b41fbbd5b136a25cbe39beedd4977163f7d36/e2e/dc985425-1c7c-4fd9-89fd-e6724d4c27d2.md", "", "", "dc985425-1c7c-4fd9-89fd-e6724d4c27d2.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/018e5cf3694a27eed9e0a56c9b9439f952746b5a/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/mt/dc985425-1c7c-4fd9-89fd-e6724d4c27d2.562879b4d3044da5adb5f238ed954a3eb20d8851.de-de.xlf", "", "", "dc985425-1c7c-4fd9-89fd-e6724d4c27d2.562879b4d3044da5adb5f238ed954a3eb20d8851.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/d03787431f5a3d598fef0f50f326d128bb146a3d/e2e/e988cb26-8b22-4b05-b024-22b86466f393.md", "", "", "e988cb26-8b22-4b05-b024-22b86466f393.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/d03787431f5a3d598fef0f50f326d128bb146a3d/e2e/e988cb26-8b22-4b05-b024-22b86466f393.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/49f891363dce45ef12fb8a60c2bc2d061cbeb1c4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/e988cb26-8b22-4b05-b024-22b86466f393.217cbbfe5e36df50d96a70fd14c971ac31104000.de-de.xlf", "", "", "e988cb26-8b22-4b05-b024-22b86466f393.217cbbfe5e36df50d96a70fd14c971ac31104000.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/243c54956fc46f4693e9a4921faacc08b9574e2d/e2e/e988cb26-8b22-4b05-b024-22b86466f393.md", "", "", "e988cb26-8b22-4b05-b024-22b86466f393.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2a2c62b0eb4ee1edf095955e40781832d0c9b312/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e988cb26-8b22-4b05-b024-22b86466f393.217cbbfe5e36df50d96a70fd14c971ac31104000.de-de.xlf", "", "", "e988cb26-8b22-4b05-b024-22b86466f393.217cbbfe5e36df50d96a70fd14c971ac31104000.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/3262cec2e5534eb1fb51addfd3cfed290b4fcb0e/e2e/05bdf957-5021-4c9f-bbc8-fe4e026e9c96.md", "", "", "a8fe6dd2-032f-4017-b24f-9e9414829f06.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B6"), "https://github.com/OpenLocalizationTest/oltest/blob/3262cec2e5534eb1fb51addfd3cfed290b4fcb0e/e2e/05bdf957-5021-4c9f-bbc8-fe4e026e9c96.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fb2c900e50bbe1255648e9682321bd740581405b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/05bdf957-5021-4c9f-bbc8-fe4e026e9c96.0d47ba9f2dce54aba0760192d3c7c3ca0bb473a8.de-de.xlf", "", "", "a8fe6dd2-032f-4017-b24f-9e9414829f06.4b0e106ca6a6d67c239a240617393aacd849fcaa.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/5ec76fc78b7923d3c111025329fa9d307811021d/e2e/73564754-f334-4035-a774-30b6bcac73ee.md", "", "", "05bdf957-5021-4c9f-bbc8-fe4e026e9c96.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B7"), "https://github.com/OpenLocalizationTest/oltest/blob/5ec76fc78b7923d3c111025329fa9d307811021d/e2e/73564754-f334-4035-a774-30b6bcac73ee.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5c0affbe02d6842e482f38043da5d48ac8a240c0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/73564754-f334-4035-a774-30b6bcac73ee.a042b5a7eea815cab0eb25d01f00914718577d55.de-de.xlf", "", "", "05bdf957-5021-4c9f-bbc8-fe4e026e9c96.0d47ba9f2dce54aba0760192d3c7c3ca0bb473a8.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/9b6d09d209c958f904d518b6d5614b54e38cf1fd/e2e/a8fe6dd2-032f-4017-b24f-9e9414829f06.md", "", "", "73564754-f334-4035-a774-30b6bcac73ee.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B8"), "https://github.com/OpenLocalizationTest/oltest/blob/9b6d09d209c958f904d518b6d5614b54e38cf1fd/e2e/a8fe6dd2-032f-4017-b24f-9e9414829f06.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D8"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5c0affbe02d6842e482f38043da5d48ac8a240c0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a8fe6dd2-032f-4017-b24f-9e9414829f06.4b0e106ca6a6d67c239a240617393aacd849fcaa.de-de.xlf", "", "", "73564754-f334-4035-a774-30b6bcac73ee.a042b5a7eea815cab0eb25d01f00914718577d55.de-de.xlf") | Out-Null
